# UserGroup.xlsx - add a "Dummy" field to the Data sheet (columns C/D) and
# document it as a new row in the Definition sheet.

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Data")
$wsDefinition = $wb.Worksheets.Item("Definition")

# --- Data sheet: add Dummy / Dummy columns (C, D) ------------------------
$wsData.Range("C1").Value = "Dummy"
$wsData.Range("D1").Value = "Dummy"
$wsData.Range("C1:D1").VerticalAlignment = "Center"

$wsData.Range("C2").Value = 1
$wsData.Range("D2").Value = 2

$wsData.Range("C3").Value = 3
$wsData.Range("D3").Value = 4

$wsData.Range("D5").Select()

# --- Definition sheet: describe the new "Dummy" field --------------------
$wsDefinition.Range("A4").Value = "Dummy"
$wsDefinition.Range("B4").Value = "int32"
$wsDefinition.Range("C4").Value = 2
$wsDefinition.Range("D4").Value = $true

$wsDefinition.Range("D7").Select()
